$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.738.63"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").Value = "1.732.00"

$cell = $ws.Range("D4")
$cell.Value = "'0.9974"
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.Value = "'242.45"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "

$cell = $ws.Range("D6")
$cell.Value = "'0.9977"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.27%  "

$cell = $ws.Range("D7")
$cell.Value = "'0.4929"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +1.13%  "

$cell = $ws.Range("D8")
$cell.Value = "'0.2623"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.20%  "

$cell = $ws.Range("D9")
$cell.Value = "'0.06221"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.32%  "

$ws.Range("D10").Value = "1.727.34"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("E11").Value = "  +2.64%  "

$cell = $ws.Range("D12")
$cell.Value = "'0.06995"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.32%  "

$cell = $ws.Range("D13")
$cell.Value = "'0.6124"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +2.38%  "

$cell = $ws.Range("D14")
$cell.Value = "'4.500"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.28%  "

$cell = $ws.Range("D15")
$cell.Value = "'77.34"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.03%  "

$cell = $ws.Range("D16")
$cell.Value = "'0.9976"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.30%  "

$ws.Range("D17").Value = "26.529.77"
$ws.Range("E17").Value = "  +0.33%  "

$cell = $ws.Range("D18")
$cell.Value = "'0.9969"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.37%  "

$cell = $ws.Range("D19")
$cell.Value = "'0.000007243"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +1.16%  "

$cell = $ws.Range("D20")
$cell.Value = "'11.44"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").Value = "1.948.16"
$ws.Range("E21").Value = "  -0.26%  "

$cell = $ws.Range("D22")
$cell.Value = "'4.490"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

$cell = $ws.Range("D23")
$cell.Value = "'8.571"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.32%  "

$cell = $ws.Range("D24")
$cell.Value = "'5.102"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.53%  "

$cell = $ws.Range("D25")
$cell.Value = "'138.06"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.54%  "

$cell = $ws.Range("D26")
$cell.Value = "'15.35"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "

$cell = $ws.Range("D27")
$cell.Value = "'1.777"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +3.16%  "

$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell = $ws.Range("D28")
$cell.Value = "'106.62"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.31%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell = $ws.Range("D29")
$cell.Value = "'1.381"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.02%  "

$cell = $ws.Range("D30")
$cell.Value = "'3.929"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.83%  "

$cell = $ws.Range("D31")
$cell.Value = "'0.07991"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "

$cell = $ws.Range("D32")
$cell.Value = "'3.676"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.32%  "

$cell = $ws.Range("D33")
$cell.Value = "'0.04483"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.89%  "

$cell = $ws.Range("D34")
$cell.Value = "'2.609"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.23%  "

$cell = $ws.Range("D35")
$cell.Value = "'1.002"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.21%  "

$cell = $ws.Range("D36")
$cell.Value = "'0.6238"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.01%  "

$cell = $ws.Range("D37")
$cell.Value = "'0.9336"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.80%  "

$cell = $ws.Range("D38")
$cell.Value = "'2.054"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +3.59%  "

$cell = $ws.Range("D39")
$cell.Value = "'2.418"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +0.48%  "

$cell = $ws.Range("D40")
$cell.Value = "'0.9973"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -0.37%  "

$ws.Range("E41").Value = "  +1.68%  "

$cell = $ws.Range("D42")
$cell.Value = "'5.615"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +3.76%  "

$cell = $ws.Range("D43")
$cell.Value = "'99.35"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.99%  "

$cell = $ws.Range("D44")
$cell.Value = "'0.3859"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.23%  "

$cell = $ws.Range("D45")
$cell.Value = "'6.910"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +3.37%  "

$cell = $ws.Range("D46")
$cell.Value = "'0.1162"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.41%  "

$cell = $ws.Range("D47")
$cell.Value = "'0.05382"
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.Value = "'7.847"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "

$cell = $ws.Range("D49")
$cell.Value = "'30.28"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "

$cell = $ws.Range("D50")
$cell.Value = "'51.78"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +1.42%  "

$cell = $ws.Range("D51")
$cell.Value = "'1.236"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -1.59%  "
